# Update "想去人数" (column F) values on the "展览" and "全部类型" sheets
# to reflect refreshed counts from the latest data pull.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (rows keyed by row number in that sheet) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 643
$ws1.Range("F4").Value  = 208
$ws1.Range("F5").Value  = 2
$ws1.Range("F6").Value  = 9755
$ws1.Range("F7").Value  = 881
$ws1.Range("F9").Value  = 1235
$ws1.Range("F10").Value = 3337
$ws1.Range("F12").Value = 110
$ws1.Range("F13").Value = 33
$ws1.Range("F14").Value = 29
$ws1.Range("F15").Value = 281
$ws1.Range("F16").Value = 522
$ws1.Range("F19").Value = 1415

# --- Sheet "全部类型" (same events, offset by one extra row) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value  = 643
$ws4.Range("F5").Value  = 208
$ws4.Range("F6").Value  = 2
$ws4.Range("F7").Value  = 9755
$ws4.Range("F8").Value  = 881
$ws4.Range("F10").Value = 1235
$ws4.Range("F11").Value = 3337
$ws4.Range("F13").Value = 110
$ws4.Range("F14").Value = 33
$ws4.Range("F15").Value = 29
$ws4.Range("F16").Value = 281
$ws4.Range("F17").Value = 522
$ws4.Range("F20").Value = 1415
